# Edit slide 2 ("Regional Drought Changes"):
#  - remove the first chart picture ("Picture 9") together with the eight
#    "*" callout rectangles that were labelling it (ids 12,13,14,15,17,18,19,20)
#  - re-center what's left (Title, "Picture 7" and its eight "*" callouts)
#    by shifting them all by the same vector, since the slide is now only
#    showing the single remaining chart.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Ids of the shapes that must be removed entirely.
$idsToDelete = @(10, 12, 13, 17, 18, 19, 20, 14, 15)

# Ids of the shapes that survive, but need to be shifted.
$idsToShift = @(21, 8, 22, 23, 24, 25, 27, 29, 32, 33)

# EMU offsets from the diff, converted to points (1 pt = 12700 EMU).
$dxPts = 2362954 / 12700
$dyPts = 144855 / 12700

# Delete first (walk backwards so indices stay valid while removing).
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($idsToDelete -contains $sh.Id) {
        $sh.Delete() | Out-Null
    }
}

# Now shift the survivors.
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($idsToShift -contains $sh.Id) {
        $sh.Left = $sh.Left + $dxPts
        $sh.Top = $sh.Top + $dyPts
    }
}
